# Commit: "add the NA's under duplicate_image_filename"
#
# Column E on the sheet is the "duplicate_image_filename" field (header in
# E1). Rows 2-21 hold the practice (p1-p4) and real stimulus rows, and all
# of them were missing a value in that column. Fill them in with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E21").Value = "NA"
